# Updated symbol list with latest price/volume data (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are forced as literal text (matching the sheet's inlineStr/text-cell
# convention) via a leading apostrophe, then the style is reset to "Normal" so
# no stray quote-prefix / number-format flag gets baked into the cell style.

$ws.Range("D2").Value = "'328.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.48%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.81%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.304"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-4.40%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08369"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.71%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.952"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-4.73%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9716"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.02%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-3.06%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1136"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.19%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1903"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.00%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09682"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-3.85%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04636"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.05%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1061"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.02%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001294"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.00%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005901"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.04%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.364"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.04%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.440"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.11%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.3358"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.43%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'8.397"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-18.28%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1354"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.22%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2725"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'9.41%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04175"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.59%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001240"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-4.48%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004436"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.62%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001300"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.69%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0002986"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-20.15%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02694"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-4.67%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05619"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.56%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007782"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.96%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1411"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.22%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007377"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-2.27%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002121"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.61%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008701"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'8.47%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3509"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006915"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.18%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.35%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003499"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'0.04%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003538"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'40.18%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.35%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.35%"
$ws.Range("E51").Style = "Normal"
